# Color the whole "Product Service - CRUD de productos + metricas con
# Actuator." bullet paragraph red (EE0000), matching the commit that adds
# w:color EE0000 to the paragraph mark run properties (w:pPr/w:rPr) and to
# every run's w:rPr in that paragraph.

$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Product*" -and $t -like "*Service*" -and $t -like "*Actuator*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Paragraph.Range includes the trailing paragraph mark, so setting the
    # font color on it also stamps w:pPr/w:rPr/w:color (paragraph-mark run
    # properties) in addition to every run in the paragraph.
    $target.Range.Font.Color = 238
    Write-Output ("Recolored paragraph: " + $target.Range.Text)
} else {
    Write-Output "Target paragraph not found."
}
